# Adds support for "min deductible, no limit" calc rules (rows 87-89)
# to the calc_rules sheet, mirroring the existing "no limit" family of
# rows already present (e.g. rows 84-86).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calc_rules")

# --- New row 87: deductible with min deductible, no limit -----------------
$ws.Cells.Item(87, 1).Value = "deductible with min deductible, no limit"
$ws.Cells.Item(87, 2).Value = 26
$ws.Cells.Item(87, 3).Value = 1
$ws.Cells.Item(87, 4).Value = 1
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Formula = '="("&C87&", "&D87&", "&E87&", "&F87&", "&G87&", "&H87&", "&I87&", "&J87&", "&K87&", "&L87&")"'

# --- New row 88: deductible with min deductible, no limit % loss ----------
$ws.Cells.Item(88, 1).Value = "deductible with min deductible, no limit % loss"
$ws.Cells.Item(88, 2).Value = 26
$ws.Cells.Item(88, 3).Value = 1
$ws.Cells.Item(88, 4).Value = 1
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 1
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Formula = '="("&C88&", "&D88&", "&E88&", "&F88&", "&G88&", "&H88&", "&I88&", "&J88&", "&K88&", "&L88&")"'

# --- New row 89: deductible with min deductible, no limit % TIV -----------
$ws.Cells.Item(89, 1).Value = "deductible with min deductible, no limit % TIV"
$ws.Cells.Item(89, 2).Value = 26
$ws.Cells.Item(89, 3).Value = 1
$ws.Cells.Item(89, 4).Value = 1
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 2
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Formula = '="("&C89&", "&D89&", "&E89&", "&F89&", "&G89&", "&H89&", "&I89&", "&J89&", "&K89&", "&L89&")"'

# --- Refresh the view: freeze panes at B56 and select the new full range --
$ws.Range("B56").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1:M89").Select()
